$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so numeric-looking values
# (e.g. "322.61", "0.535") are not auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "47.133.03"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "2.489.95"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "322.61"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "108.54"
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "38.95"
$ws.Range("E10").Value = "  +7.33%  "
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "18.40"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "7.19"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "2.882.52"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "2.496.13"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "0.851"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "47.061.41"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "6.60"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "0.0₃0937"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "2.70"
$ws.Range("E22").Value = "  +12.41%  "
$ws.Range("D23").Value = "70.69"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "247.14"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "25.87"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "10.05"
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("D30").Value = "0.141"
$ws.Range("E30").Value = "  +9.47%  "
$ws.Range("D31").Value = "35.17"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").Value = "49.92"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").Value = "19.98"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").Value = "1.97"
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("D38").Value = "4.69"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("D39").Value = "2.98"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").Value = "120.82"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D43").Value = "21.29"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").Value = "0.0296"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").Value = "1.992.99"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "3.04"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.79"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "9.11"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "5.18"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "56.40"
$ws.Range("E51").Value = "  +2.69%  "
